$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the Fitness column (C) values for rows 2 through 79 (Generation 0-77)
# to the corrected constant value of 7293, per "correction in sa algorithm".
$ws.Range("C2:C79").Value = 7293
